# RESO APPOINT 1ST COSEC.docx - merge split (proofErr-fragmented) runs back
# into single runs, and add Jinja-style {% for %} / {% endfor %} loop markers
# around the directors table.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Body: "...Directors of {{ company_name }} (the Company)..."
#    Collapse the proofErr-split "{{ ", "company", "_name", " }}" runs into
#    one run with the full placeholder text.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("{{ company_name }}", $false, $false, $false, $false, $false, `
    $true, 1, $false, "{{ company_name }}", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Body: "Name: Mohammad Syafuan Bin Nordin" - merge the name runs.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Mohammad Syafuan Bin Nordin", $false, $false, $false, $false, $false, `
    $true, 1, $false, "Mohammad Syafuan Bin Nordin", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) Body: "Firm Name: AMR Secretarial Services Sdn. Bhd." - merge the firm
#    name runs.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("AMR Secretarial Services Sdn. Bhd.", $false, $false, $false, $false, $false, `
    $true, 1, $false, "AMR Secretarial Services Sdn. Bhd.", 2) | Out-Null

# ---------------------------------------------------------------------------
# 4) Body: "Date: {{ incorporation_date }}" - merge the date placeholder
#    runs.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("{{ incorporation_date }}", $false, $false, $false, $false, $false, `
    $true, 1, $false, "{{ incorporation_date }}", 2) | Out-Null

# ---------------------------------------------------------------------------
# 5) Table cell (row 1, left column): add the "if row.left" condition to the
#    two placeholders (text content actually changes here, not just a run
#    merge).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("{{ row.left.line }}", $false, $false, $false, $false, $false, `
    $true, 1, $false, "{{ row.left.line if row.left }}", 2) | Out-Null

$d.Content.Find.Execute("{{ row.left.name }}", $false, $false, $false, $false, $false, `
    $true, 1, $false, "{{ row.left.name if row.left }}", 2) | Out-Null

# ---------------------------------------------------------------------------
# 6) Table cell (row 1, right column): merge the already-conditioned
#    placeholders' runs (no text change).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("{{ row.right.line if row.right }}", $false, $false, $false, $false, $false, `
    $true, 1, $false, "{{ row.right.line if row.right }}", 2) | Out-Null

$d.Content.Find.Execute("{{ row.right.name if row.right }}", $false, $false, $false, $false, $false, `
    $true, 1, $false, "{{ row.right.name if row.right }}", 2) | Out-Null

# ---------------------------------------------------------------------------
# 7) Header: "{{ company_name }}" and "{{ ssm_number }}" - merge runs.
# ---------------------------------------------------------------------------
$sec = $d.Sections(1)
$hdr = $sec.Headers(1)
$hdr.Range.Find.Execute("{{ company_name }}", $false, $false, $false, $false, $false, `
    $true, 1, $false, "{{ company_name }}", 2) | Out-Null
$hdr.Range.Find.Execute("{{ ssm_number }}", $false, $false, $false, $false, $false, `
    $true, 1, $false, "{{ ssm_number }}", 2) | Out-Null

# ---------------------------------------------------------------------------
# 8) Insert a new paragraph containing "{% for row in director_rows %}"
#    directly above the directors table (after the existing blank centered
#    paragraph that precedes the table).
#
#    NOTE: Table.Range.Paragraphs(...) is unreliable near table boundaries in
#    this runtime (it can report stale/out-of-range positions), so the
#    paragraph immediately before/after the table is located by scanning the
#    document's Paragraphs collection for a Range whose Start/End abuts the
#    table's Range - and, to dodge a degenerate (zero-length) phantom
#    paragraph that appears right after the table, the "after" search also
#    requires a non-empty span.
# ---------------------------------------------------------------------------
$tbl = $d.Tables(1)

$beforePara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.End -eq $tbl.Range.Start) {
        $beforePara = $p
    }
}
$beforePara.Range.InsertParagraphAfter() | Out-Null

$forPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.End -eq $tbl.Range.Start) {
        $forPara = $p
    }
}
$forPara.Range.ParagraphFormat.Alignment = 0
$forPara.Range.Text = "{% for row in director_rows %}"
$forRange = $forPara.Range
$forRange.Font.NameAscii = "Times New Roman"
$forRange.Font.NameFarEast = "Times New Roman"
$forRange.Font.NameOther = "Times New Roman"
$forRange.Font.NameBi = "Times New Roman"
$forRange.Font.Bold = $true
$forRange.Font.BoldBi = $true

# ---------------------------------------------------------------------------
# 9) Add a "{% endfor %}" run to the (currently empty) paragraph that sits
#    immediately after the directors table.
# ---------------------------------------------------------------------------
$afterPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Start -eq $tbl.Range.End -and $p.Range.End -gt $p.Range.Start) {
        $afterPara = $p
    }
}
$endRange = $afterPara.Range.Duplicate
$endRange.Collapse(1)
$endRange.InsertBefore("{% endfor %}")

$afterPara2 = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Start -eq $tbl.Range.End -and $p.Range.Text -eq "{% endfor %}`r") {
        $afterPara2 = $p
    }
}
$endFontRange = $afterPara2.Range
$endFontRange.Font.NameAscii = "Times New Roman"
$endFontRange.Font.NameFarEast = "Times New Roman"
$endFontRange.Font.NameOther = "Times New Roman"
$endFontRange.Font.NameBi = "Times New Roman"
$endFontRange.Font.Bold = $true
$endFontRange.Font.BoldBi = $true

Write-Host "Done"
